# Edit: insert two new Kiwi price rows (positions 352 and 353) into the
# "Vega Modelo de Temuco - Kiwi" sheet, shifting the previously existing
# rows 352-417 down to 354-419 and extending the used range to A1:T419.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 352..417 down by 2 rows, creating two blank rows
# at 352 and 353.
$ws.Rows.Item(352).Resize(2).Insert()

# Populate the two newly inserted rows with the new data.
$ws.Range("A352").Value = 10
$ws.Range("B352").Value = "Vega Modelo de Temuco"
$ws.Range("C352").Value = "La Araucanía"
$ws.Range("D352").Value = 44694
$ws.Range("E352").Value = 9
$ws.Range("F352").Value = "Fruta"
$ws.Range("G352").Value = 100101
$ws.Range("H352").Value = "Berries"
$ws.Range("I352").Value = 100101007
$ws.Range("J352").Value = "Kiwi"
$ws.Range("K352").Value = "Hayward"
$ws.Range("L352").Value = "Especial"
$ws.Range("M352").Value = 80
$ws.Range("N352").Value = 18000
$ws.Range("O352").Value = 18000
$ws.Range("P352").Value = 18000
$ws.Range("Q352").Value = "$/bandeja 18 kilos"
$ws.Range("R352").Value = "Región de O'Higgins"
$ws.Range("S352").Value = 1000
$ws.Range("T352").Value = 18

$ws.Range("A353").Value = 10
$ws.Range("B353").Value = "Vega Modelo de Temuco"
$ws.Range("C353").Value = "La Araucanía"
$ws.Range("D353").Value = 44694
$ws.Range("E353").Value = 9
$ws.Range("F353").Value = "Fruta"
$ws.Range("G353").Value = 100101
$ws.Range("H353").Value = "Berries"
$ws.Range("I353").Value = 100101007
$ws.Range("J353").Value = "Kiwi"
$ws.Range("K353").Value = "Hayward"
$ws.Range("L353").Value = "Primera"
$ws.Range("M353").Value = 95
$ws.Range("N353").Value = 15000
$ws.Range("O353").Value = 15000
$ws.Range("P353").Value = 15000
$ws.Range("Q353").Value = "$/bandeja 18 kilos"
$ws.Range("R353").Value = "Región de O'Higgins"
$ws.Range("S353").Value = 833
$ws.Range("T353").Value = 18
